$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 671, pushing existing data (rows 671-723) down to 673-725.
$ws.Rows.Item(671).Insert()
$ws.Rows.Item(671).Insert()

# New row 671: "Pintón" quality record
$ws.Cells.Item(671, 1).Value = 5
$ws.Cells.Item(671, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(671, 3).Value = "Maule"
$ws.Cells.Item(671, 4).Value = 44783
$ws.Cells.Item(671, 5).Value = 7
$ws.Cells.Item(671, 6).Value = "Fruta"
$ws.Cells.Item(671, 7).Value = 100108
$ws.Cells.Item(671, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(671, 9).Value = 100108006
$ws.Cells.Item(671, 10).Value = "Plátano"
$ws.Cells.Item(671, 11).Value = "Sin especificar"
$ws.Cells.Item(671, 12).Value = "Pintón"
$ws.Cells.Item(671, 13).Value = 800
$ws.Cells.Item(671, 14).Value = 19000
$ws.Cells.Item(671, 15).Value = 19000
$ws.Cells.Item(671, 16).Value = 19000
$ws.Cells.Item(671, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(671, 18).Value = "Ecuador"
$ws.Cells.Item(671, 19).Value = 950
$ws.Cells.Item(671, 20).Value = 20

# New row 672: "Primera Pintón" quality record
$ws.Cells.Item(672, 1).Value = 5
$ws.Cells.Item(672, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(672, 3).Value = "Maule"
$ws.Cells.Item(672, 4).Value = 44783
$ws.Cells.Item(672, 5).Value = 7
$ws.Cells.Item(672, 6).Value = "Fruta"
$ws.Cells.Item(672, 7).Value = 100108
$ws.Cells.Item(672, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(672, 9).Value = 100108006
$ws.Cells.Item(672, 10).Value = "Plátano"
$ws.Cells.Item(672, 11).Value = "Sin especificar"
$ws.Cells.Item(672, 12).Value = "Primera Pintón"
$ws.Cells.Item(672, 13).Value = 500
$ws.Cells.Item(672, 14).Value = 20000
$ws.Cells.Item(672, 15).Value = 20000
$ws.Cells.Item(672, 16).Value = 20000
$ws.Cells.Item(672, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(672, 18).Value = "Ecuador"
$ws.Cells.Item(672, 19).Value = 1000
$ws.Cells.Item(672, 20).Value = 20

Write-Host "done"
